# bug449-template.docx : turn the M2Doc field ("{ m:'doc.html'.fromHTMLURI() }")
# that is stored as a real Word field (fldChar begin/instrText.../fldChar end)
# into plain template text made of literal runs: "{", "m", ":", "'", "doc.html",
# "'.fromHTMLURI()", "}" - keeping the _GoBack bookmark exactly where it was,
# straddling the boundary between "doc.html" and "'.fromHTMLURI()".

$d = $word.ActiveDocument

# The field lives alone in the document's 2nd paragraph.
$fieldParagraph = $d.Paragraphs(2)
$paraStart = $fieldParagraph.Range.Start

# Deleting the field removes the fldChar begin/end runs, every instrText run
# and the _GoBack bookmark that sat inside the field code - the paragraph
# becomes empty (just the paragraph mark remains).
$field = $d.Fields.Item(1)
$field.Delete()

# Re-type the whole instruction as plain, visible text - this is what the
# template text now looks like once it is no longer a field: the field
# delimiters become literal "{" / "}" characters.
$plainText = "{m:'doc.html'.fromHTMLURI()}"
$insertionPoint = $d.Range($paraStart, $paraStart)
$insertionPoint.InsertAfter($plainText)

# Right now all of that text sits in a single run. Word only ever breaks a
# run in two when something non-text (a bookmark, a field, ...) is anchored
# in the middle of it, so we drop a temporary bookmark at every boundary we
# need and then discard the ones we do not actually want to keep. This
# leaves one run per segment: "{", "m", ":", "'", "doc.html",
# "'.fromHTMLURI()", "}".
$segments = @("{", "m", ":", "'", "doc.html", "'.fromHTMLURI()", "}")

$offset = 0
$boundaryNames = @()
for ($i = 0; $i -lt ($segments.Length - 1); $i++) {
    $offset = $offset + $segments[$i].Length
    # The boundary right after "doc.html" is where _GoBack used to be (and
    # needs to be again); every other boundary only exists to split runs and
    # gets removed again below.
    if ($i -eq 4) {
        $boundaryNames += "_GoBack"
    } else {
        $boundaryNames += "TmpSplit$i"
    }
    $pos = $paraStart + $offset
    $boundaryRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add($boundaryNames[$i], $boundaryRange) | Out-Null
}

foreach ($name in $boundaryNames) {
    if ($name -ne "_GoBack") {
        $d.Bookmarks($name).Delete()
    }
}

Write-Output "Field replaced with plain text: [$($fieldParagraph.Range.Text)]"
